$d = $word.ActiveDocument

# 1) Merge the "As a data scientist... develop my data science skills." runs
#    (also drops the spell-check proofErr wrapper around "skills").
$r = $d.Content
$r.Find.Execute(
    "As a data scientist, I want to read about models used in the project, so that I can develop my data science skills.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a data scientist, I want to read about models used in the project, so that I can develop my data science skills.",
    2)

# 2) Merge the "As a data scientist, I want to choose a model..." + "." runs
$r2 = $d.Content
$r2.Find.Execute(
    "As a data scientist, I want to choose a model, so that I can compare the results of different clustering algorithms.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a data scientist, I want to choose a model, so that I can compare the results of different clustering algorithms.",
    2)

# 3) Merge the "# " + "Generated charts (report)..." runs
$r3 = $d.Content
$r3.Find.Execute(
    "# Generated charts (report) will be interactive, enabling the user to perform zooming, hovering, and more to display preferred information.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "# Generated charts (report) will be interactive, enabling the user to perform zooming, hovering, and more to display preferred information.",
    2)
